# Auto-generated edit script applying numeric corrections to Sheets/Jenova_Profits.xlsx
# (multi-sheet workbook: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 76: Warding Off Temptation / Enchanted Hardsilver Ink
$ws.Range("H76").Value = 66674550
$ws.Range("I76").Value = 7660
$ws.Range("K76").Value = 7660
$ws.Range("M76").Value = -7345

# Row 79: The Garden of Arcane Delights (L) / Enchanted Hardsilver Ink
$ws.Range("H79").Value = 66674550
$ws.Range("I79").Value = 7660
$ws.Range("K79").Value = 7660
$ws.Range("M79").Value = -6568

# Row 132: Fast-forwarding Flora / Growth Formula Lambda
$ws.Range("H132").Value = 2845.9355
$ws.Range("I132").Value = 2934.2593
$ws.Range("K132").Value = 8802.777900000001
$ws.Range("M132").Value = -6272.777900000001

# Row 141: Remedy for Reason / Grade 1 Gemdraught of Mind
$ws.Range("H141").Value = 8374.25
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
# Row 61: Dealing with the Tough Stuff / Cobalt Ingot
$ws.Range("H61").Value = 2712.8838
$ws.Range("I61").Value = 2829.2727
$ws.Range("J61").Value = 2328.8
$ws.Range("K61").Value = 2829.2727
$ws.Range("L61").Value = 2328.8
$ws.Range("M61").Value = -2617.2727
$ws.Range("N61").Value = -2752.8

# Row 74: As the Bolt Flies / Titanium Nugget
$ws.Range("H74").Value = 2489.3157
$ws.Range("I74").Value = 2802.9167
$ws.Range("K74").Value = 2802.9167
$ws.Range("M74").Value = -1928.9167

# Row 77: Heavy Metal Banned (L) / Titanium Nugget
$ws.Range("H77").Value = 2489.3157
$ws.Range("I77").Value = 2802.9167
$ws.Range("K77").Value = 14014.5835
$ws.Range("M77").Value = -9646.583500000001

# Row 122: Haste for High Durium / High Durium Nugget
$ws.Range("H122").Value = 2918.66
$ws.Range("I122").Value = 1354.0385
$ws.Range("J122").Value = 4613.6665
$ws.Range("K122").Value = 4062.1155
$ws.Range("L122").Value = 13840.9995
$ws.Range("M122").Value = -1612.1155
$ws.Range("N122").Value = -18740.9995

# Row 125: The Incomplete Costume / High Durium Armor of Fending
$ws.Range("H125").Value = 88945.75
$ws.Range("J125").Value = 88945.75
$ws.Range("L125").Value = 88945.75
$ws.Range("N125").Value = -98785.75

# Row 136: Metal with Mettle / Cobalt Tungsten Ingot
$ws.Range("H136").Value = 2712.8838
$ws.Range("I136").Value = 2829.2727
$ws.Range("J136").Value = 2328.8
$ws.Range("K136").Value = 8487.8181
$ws.Range("L136").Value = 6986.400000000001
$ws.Range("M136").Value = -5937.8181
$ws.Range("N136").Value = -12086.4

$ws = $wb.Worksheets.Item("BSM")
# Row 81: Diamond Sawdust / Titanium Battleaxe
$ws.Range("H81").Value = 61593
$ws.Range("J81").Value = 79889.5
$ws.Range("L81").Value = 79889.5
$ws.Range("N81").Value = -82011.5

# Row 84: I'm a Lumberjack and I'm Okay (L) / Titanium Battleaxe
$ws.Range("H84").Value = 61593
$ws.Range("J84").Value = 79889.5
$ws.Range("L84").Value = 239668.5
$ws.Range("N84").Value = -250276.5

# Row 107: The Gold Experience / Deepgold Nugget
$ws.Range("H107").Value = 911396.75
$ws.Range("I107").Value = 1793.875
$ws.Range("J107").Value = 3337004.2
$ws.Range("K107").Value = 1793.875
$ws.Range("L107").Value = 3337004.2
$ws.Range("M107").Value = 126.125
$ws.Range("N107").Value = -3340844.2

# Row 134: Ruthenium Supremium / Ruthenium Ingot
$ws.Range("H134").Value = 20824.754
$ws.Range("I134").Value = 3121.0212
$ws.Range("K134").Value = 9363.063600000001
$ws.Range("M134").Value = -6828.063600000001

$ws = $wb.Worksheets.Item("CRP")
# Row 4: A Clogful of Camaraderie / Maple Clogs
$ws.Range("H4").Value = 11431715

# Row 58: You Do the Heavy Lifting / Mahogany Lumber
$ws.Range("H58").Value = 2166.375
$ws.Range("I58").Value = 2762.0667
$ws.Range("K58").Value = 2762.0667
$ws.Range("M58").Value = -2559.0667

# Row 122: Timber of Tenkonto / Horse Chestnut Lumber
$ws.Range("H122").Value = 2459
$ws.Range("I122").Value = 2374.7144
$ws.Range("K122").Value = 7124.1432
$ws.Range("M122").Value = -4674.1432

# Row 136: Turali Quality / Dark Mahogany Lumber
$ws.Range("H136").Value = 2166.375
$ws.Range("I136").Value = 2762.0667
$ws.Range("K136").Value = 8286.2001
$ws.Range("M136").Value = -5736.2001

$ws = $wb.Worksheets.Item("CUL")
# Row 50: Moving Up in the World / Rolanberry Cheese
$ws.Range("H50").Value = 38711.383
$ws.Range("I50").Value = 435.25
$ws.Range("J50").Value = 55723
$ws.Range("K50").Value = 1305.75
$ws.Range("L50").Value = 167169
$ws.Range("M50").Value = -824.75
$ws.Range("N50").Value = -168131

# Row 53: Rolanberry Fields Forever / Rolanberry Cheese
$ws.Range("H53").Value = 38711.383
$ws.Range("I53").Value = 435.25
$ws.Range("J53").Value = 55723
$ws.Range("K53").Value = 1305.75
$ws.Range("L53").Value = 167169
$ws.Range("M53").Value = -824.75
$ws.Range("N53").Value = -168131

# Row 87: Soup That Eats Like a Knight / Clam Chowder
$ws.Range("H87").Value = 13311.889
$ws.Range("J87").Value = 16209.632
$ws.Range("L87").Value = 48628.896
$ws.Range("N87").Value = -51124.896

# Row 90: Like Ma Used to Make (L) / Clam Chowder
$ws.Range("H90").Value = 13311.889
$ws.Range("J90").Value = 16209.632
$ws.Range("L90").Value = 145886.688
$ws.Range("N90").Value = -158366.688

# Row 107: Slippery Service / Frantoio Oil
$ws.Range("H107").Value = 2990.577
$ws.Range("J107").Value = 4191.4375
$ws.Range("L107").Value = 12574.3125
$ws.Range("N107").Value = -16414.3125

# Row 122: Salt of the North / Northern Sea Salt
$ws.Range("H122").Value = 2225.3
$ws.Range("I122").Value = 1382.4445
$ws.Range("K122").Value = 12442.0005
$ws.Range("M122").Value = -9992.0005

# Row 131: The Mountain Steeped / Tsai tou Vounou
$ws.Range("H131").Value = 2783.889
$ws.Range("I131").Value = 1588.3334
$ws.Range("J131").Value = 3182.4075
$ws.Range("K131").Value = 4765.0002
$ws.Range("L131").Value = 9547.2225
$ws.Range("M131").Value = 274.9997999999996
$ws.Range("N131").Value = -19627.2225

# Row 132: More Mezcal / Cooking Mezcal
$ws.Range("H132").Value = 481527.3
$ws.Range("J132").Value = 628895.5600000001
$ws.Range("L132").Value = 5660060.040000001
$ws.Range("N132").Value = -5665120.040000001

$ws = $wb.Worksheets.Item("GSM")
# Row 102: Put the Metal to the Peddle / Durium Ingot
$ws.Range("H102").Value = 2099.4285
$ws.Range("I102").Value = 1095.138
$ws.Range("K102").Value = 1095.138
$ws.Range("M102").Value = 526.8620000000001

# Row 122: Awarding Academic Excellence / Ametrine
$ws.Range("H122").Value = 2958.5862
$ws.Range("I122").Value = 2047.9565
$ws.Range("J122").Value = 6449.3335
$ws.Range("K122").Value = 6143.8695
$ws.Range("L122").Value = 19348.0005
$ws.Range("M122").Value = -3693.8695
$ws.Range("N122").Value = -24248.0005

# Row 123: Workplace Workout / Ametrine Ring of Fending
$ws.Range("H123").Value = 28332.5
$ws.Range("J123").Value = 28332.5
$ws.Range("L123").Value = 28332.5
$ws.Range("N123").Value = -33232.5

# Row 126: Gold Rush Order / Phrygian Gold Ingot
$ws.Range("H126").Value = 4295.8667
$ws.Range("I126").Value = 4394.4
$ws.Range("J126").Value = 4246.6
$ws.Range("K126").Value = 13183.2
$ws.Range("L126").Value = 12739.8
$ws.Range("M126").Value = -10713.2
$ws.Range("N126").Value = -17679.8

# Row 136: Shiny and Good / Pink Beryl
$ws.Range("H136").Value = 32924.285
$ws.Range("J136").Value = 32924.285
$ws.Range("L136").Value = 98772.85500000001
$ws.Range("N136").Value = -103872.855

$ws = $wb.Worksheets.Item("LTW")
# Row 2: Red in the Head / Leather Calot
$ws.Range("H2").Value = 66670830
$ws.Range("I2").Value = 200002500
$ws.Range("K2").Value = 200002500
$ws.Range("M2").Value = -200002388

# Row 7: Tan Before the Ban / Leather
$ws.Range("H7").Value = 6570.28
$ws.Range("I7").Value = 6257.65
$ws.Range("K7").Value = 6257.65
$ws.Range("M7").Value = -6145.65

# Row 16: Saddle Sore / Hard Leather
$ws.Range("H16").Value = 531.6667
$ws.Range("I16").Value = 489.0909
$ws.Range("J16").Value = 1000
$ws.Range("K16").Value = 489.0909
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = -319.0909
$ws.Range("N16").Value = -1340

# Row 40: Best Served Toad / Toad Leather
$ws.Range("H40").Value = 4852.3335
$ws.Range("I40").Value = 2540.818
$ws.Range("K40").Value = 2540.818
$ws.Range("M40").Value = -2404.818

# Row 82: Trainin' the Neck / Dragon Leather
$ws.Range("H82").Value = 1541.8125
$ws.Range("I82").Value = 1244.6364
$ws.Range("J82").Value = 2195.6
$ws.Range("K82").Value = 1244.6364
$ws.Range("L82").Value = 2195.6
$ws.Range("M82").Value = -883.6364000000001
$ws.Range("N82").Value = -2917.6

# Row 85: Training Is Only Skintight (L) / Dragon Leather
$ws.Range("H85").Value = 1541.8125
$ws.Range("I85").Value = 1244.6364
$ws.Range("J85").Value = 2195.6
$ws.Range("K85").Value = 1244.6364
$ws.Range("L85").Value = 2195.6
$ws.Range("M85").Value = 3.363599999999906
$ws.Range("N85").Value = -4691.6

# Row 122: Hell on Leather / Gaja Leather
$ws.Range("H122").Value = 4610.963
$ws.Range("I122").Value = 4285.091
$ws.Range("J122").Value = 5123.048
$ws.Range("K122").Value = 12855.273
$ws.Range("L122").Value = 15369.144
$ws.Range("M122").Value = -10405.273
$ws.Range("N122").Value = -20269.144

# Row 126: Battered Books / Saiga Leather
$ws.Range("H126").Value = 6570.28
$ws.Range("I126").Value = 6257.65
$ws.Range("K126").Value = 18772.95
$ws.Range("M126").Value = -16302.95

# Row 132: Tenets of Tanning / Silver Lobo Leather
$ws.Range("H132").Value = 3039.682
$ws.Range("I132").Value = 2601
$ws.Range("J132").Value = 4531.2
$ws.Range("K132").Value = 7803
$ws.Range("L132").Value = 13593.6
$ws.Range("M132").Value = -5273
$ws.Range("N132").Value = -18653.6

# Row 136: Respect for Br'aax / Br'aax Leather
$ws.Range("H136").Value = 169438.4
$ws.Range("I136").Value = 290176.8
$ws.Range("K136").Value = 870530.3999999999
$ws.Range("M136").Value = -867980.3999999999

$ws = $wb.Worksheets.Item("WVR")
# Row 107: Flax Wax / Bright Linen Yarn
$ws.Range("H107").Value = 1473.381
$ws.Range("I107").Value = 1668.5
$ws.Range("K107").Value = 5005.5
$ws.Range("M107").Value = -3085.5

# Row 116: All-purpose Overgarments / Ovim Wool Tunic of Aiming
$ws.Range("H116").Value = 94832
$ws.Range("J116").Value = 94832
$ws.Range("L116").Value = 94832
$ws.Range("N116").Value = -104010

# Row 122: Heavy Armoire / Dark Hempen Cloth
$ws.Range("H122").Value = 25643778
$ws.Range("I122").Value = 37039396
$ws.Range("K122").Value = 111118188
$ws.Range("M122").Value = -111115738

# Row 126: A Polished Purchase / Snow Linen
$ws.Range("H126").Value = 1127.7693
$ws.Range("I126").Value = 1127.7693
$ws.Range("K126").Value = 3383.3079
$ws.Range("M126").Value = -913.3078999999998

# Row 136: Weaving the Envelope / Sarcenet Cloth
$ws.Range("H136").Value = 40317.938
$ws.Range("I136").Value = 9764.849
$ws.Range("K136").Value = 29294.547
$ws.Range("M136").Value = -26744.547
